# overview.pptx edit: "generator.yml" -> "genlayer.yml" (#482)
# plus the handout/notes master "datetimeFigureOut" cached date text
# bumping from 13/09/2022 to 14/09/2022.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2, shape "Flowchart: Multidocument 11": the label reading
#    "*.generator.yml / Defines generated / project parts" becomes
#    "*.genlayer.yml / Defines generated / project parts".
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp = $s2.Shapes.Item(29)
$tr = $shp.TextFrame.TextRange

# "*.generator.yml..." -> replace the "generator." portion (chars 3-12,
# right after the literal "*.") with "genlayer." ; this keeps the
# leading "*." run untouched and leaves the trailing "yml" as its own
# run, matching how the text is split in the authored deck.
$editRange = $tr.Characters(3, 10)
if ($editRange.Text -eq "generator.") {
    $editRange.Text = "genlayer."
}

# ---------------------------------------------------------------------
# 2) Handout master + notes master "Date Placeholder" fields: cached
#    datetimeFigureOut text 13/09/2022 -> 14/09/2022.
# ---------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    $hmDate = $hm.Shapes.Item(2)
    if ($hmDate.TextFrame.TextRange.Text -eq "13/09/2022") {
        $hmDate.TextFrame.TextRange.Text = "14/09/2022"
    }
} catch {
    Write-Host "handout master date field could not be edited: " $_.Exception.Message
}

try {
    $nm = $p.NotesMaster
    $nmDate = $nm.Shapes.Item(2)
    if ($nmDate.TextFrame.TextRange.Text -eq "13/09/2022") {
        $nmDate.TextFrame.TextRange.Text = "14/09/2022"
    }
} catch {
    Write-Host "notes master date field could not be edited: " $_.Exception.Message
}
